$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value as literal text, even when the string looks like a
# number (e.g. "580.34", "1.00", "3.180.64"), matching the source data which
# stores these as plain text (inline strings), not numeric cells.
# We briefly force Text number format so Excel does not auto-convert the
# string into a numeric value (which would also destroy formatting such as
# trailing zeros), then clear the format again so the cell is left exactly
# as it was before (no lingering style / number-format changes).
function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextValue "D2" "65.093.51"
Set-TextValue "E2" "  +2.03%  "
Set-TextValue "D3" "3.177.30"
Set-TextValue "E3" "  +4.14%  "
Set-TextValue "E4" "  +0.08%  "
Set-TextValue "D5" "580.34"
Set-TextValue "E5" "  +3.95%  "
Set-TextValue "D6" "151.18"
Set-TextValue "E6" "  +6.20%  "
Set-TextValue "E7" "  -0.03%  "
Set-TextValue "D8" "3.175.06"
Set-TextValue "E8" "  +4.15%  "
Set-TextValue "E9" "  +4.06%  "
Set-TextValue "E10" "  +5.99%  "
Set-TextValue "D11" "6.21"
Set-TextValue "E11" "  +0.17%  "
Set-TextValue "D12" "0.504"
Set-TextValue "E12" "  +2.20%  "
Set-TextValue "D13" "0.0000273"
Set-TextValue "E13" "  +18.96%  "
Set-TextValue "D14" "37.92"
Set-TextValue "E14" "  +6.42%  "
Set-TextValue "D15" "3.703.80"
Set-TextValue "E15" "  +4.32%  "
Set-TextValue "D16" "65.185.49"
Set-TextValue "E16" "  +2.13%  "
Set-TextValue "D17" "3.179.84"
Set-TextValue "E17" "  +4.27%  "
Set-TextValue "D18" "7.18"
Set-TextValue "E18" "  +5.88%  "
Set-TextValue "E19" "  +1.60%  "
Set-TextValue "D20" "513.78"
Set-TextValue "E20" "  +8.15%  "
Set-TextValue "D21" "14.93"
Set-TextValue "E21" "  +6.56%  "
Set-TextValue "D22" "0.732"
Set-TextValue "E22" "  +7.12%  "
Set-TextValue "D23" "15.24"
Set-TextValue "E23" "  +4.90%  "
Set-TextValue "D24" "7.81"
Set-TextValue "E24" "  +3.51%  "
Set-TextValue "D25" "85.41"
Set-TextValue "E25" "  +3.54%  "
Set-TextValue "D26" "0.999"
Set-TextValue "E26" "  -0.01%  "
Set-TextValue "E27" "  +11.39%  "
Set-TextValue "E28" "  +4.91%  "
Set-TextValue "E29" "  +7.92%  "
Set-TextValue "D30" "27.99"
Set-TextValue "E30" "  +6.68%  "
Set-TextValue "E31" "  +14.01%  "
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.06%  "
Set-TextValue "D33" "1.20"
Set-TextValue "E33" "  +5.42%  "
Set-TextValue "D34" "6.36"
Set-TextValue "E34" "  +10.62%  "
Set-TextValue "D35" "6.61"
Set-TextValue "E35" "  +6.32%  "
Set-TextValue "D36" "55.66"
Set-TextValue "E36" "  +2.12%  "
Set-TextValue "D37" "0.0906"
Set-TextValue "E37" "  +11.54%  "
Set-TextValue "B38" "dogwifhat"
Set-TextValue "C38" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextValue "D38" "3.16"
Set-TextValue "E38" "  +12.11%  "
Set-TextValue "B39" "Bittensor"
Set-TextValue "C39" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue "D39" "475.93"
Set-TextValue "E39" "  +6.49%  "
Set-TextValue "E40" "  +3.47%  "
Set-TextValue "E41" "  +4.96%  "
Set-TextValue "D42" "3.071.68"
Set-TextValue "E42" "  +2.16%  "
Set-TextValue "D43" "0.120"
Set-TextValue "E43" "  +2.33%  "
Set-TextValue "E44" "  +7.32%  "
Set-TextValue "D45" "2.43"
Set-TextValue "E45" "  +8.09%  "
Set-TextValue "D46" "29.18"
Set-TextValue "E46" "  +4.06%  "
Set-TextValue "D47" "0.0₃0613"
Set-TextValue "E47" "  +19.54%  "
Set-TextValue "E48" "  -0.03%  "
Set-TextValue "E49" "  +2.33%  "
Set-TextValue "D50" "2.26"
Set-TextValue "E50" "  +8.48%  "
Set-TextValue "D51" "120.97"
Set-TextValue "E51" "  +2.72%  "
